$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.000.46"
Set-TextValue $ws.Range("E2") "  +1.08%  "

Set-TextValue $ws.Range("D3") "3.934.46"
Set-TextValue $ws.Range("E3") "  +0.21%  "

Set-TextValue $ws.Range("E4") "  +0.05%  "

Set-TextValue $ws.Range("D5") "487.54"
Set-TextValue $ws.Range("E5") "  -0.31%  "

Set-TextValue $ws.Range("D6") "146.34"
Set-TextValue $ws.Range("E6") "  +0.42%  "

Set-TextValue $ws.Range("D7") "0.622"
Set-TextValue $ws.Range("E7") "  -0.63%  "

Set-TextValue $ws.Range("E8") "  +0.02%  "

Set-TextValue $ws.Range("E9") "  +0.75%  "

Set-TextValue $ws.Range("E10") "  +3.75%  "

Set-TextValue $ws.Range("D12") "42.90"
Set-TextValue $ws.Range("E12") "  +0.67%  "

Set-TextValue $ws.Range("D13") "10.44"
Set-TextValue $ws.Range("E13") "  -1.85%  "

Set-TextValue $ws.Range("D14") "4.560.29"
Set-TextValue $ws.Range("E14") "  +0.15%  "

Set-TextValue $ws.Range("D15") "3.928.45"
Set-TextValue $ws.Range("E15") "  +0.44%  "

Set-TextValue $ws.Range("D16") "14.28"
Set-TextValue $ws.Range("E16") "  -4.08%  "

Set-TextValue $ws.Range("E17") "  -0.66%  "

Set-TextValue $ws.Range("D18") "19.89"
Set-TextValue $ws.Range("E18") "  -1.04%  "

Set-TextValue $ws.Range("E19") "  +2.64%  "

Set-TextValue $ws.Range("D20") "69.026.08"
Set-TextValue $ws.Range("E20") "  +0.95%  "

Set-TextValue $ws.Range("D21") "437.03"
Set-TextValue $ws.Range("E21") "  -2.18%  "

Set-TextValue $ws.Range("D22") "3.45"
Set-TextValue $ws.Range("E22") "  +1.81%  "

Set-TextValue $ws.Range("D23") "14.58"
Set-TextValue $ws.Range("E23") "  -1.52%  "

Set-TextValue $ws.Range("D24") "12.33"
Set-TextValue $ws.Range("E24") "  +14.21%  "

Set-TextValue $ws.Range("D25") "89.23"
Set-TextValue $ws.Range("E25") "  +0.49%  "

Set-TextValue $ws.Range("E26") "  +2.71%  "

Set-TextValue $ws.Range("D27") "11.10"
Set-TextValue $ws.Range("E27") "  -3.16%  "

Set-TextValue $ws.Range("D28") "37.13"
Set-TextValue $ws.Range("E28") "  -4.38%  "

Set-TextValue $ws.Range("D29") "5.65"
Set-TextValue $ws.Range("E29") "  -3.93%  "

Set-TextValue $ws.Range("D30") "710.41"
Set-TextValue $ws.Range("E30") "  +2.66%  "

Set-TextValue $ws.Range("E31") "  +1.33%  "

Set-TextValue $ws.Range("D32") "13.54"
Set-TextValue $ws.Range("E32") "  +0.53%  "

Set-TextValue $ws.Range("D33") "2.90"
Set-TextValue $ws.Range("E33") "  +1.18%  "

Set-TextValue $ws.Range("D34") "0.484"
Set-TextValue $ws.Range("E34") "  +30.68%  "

Set-TextValue $ws.Range("D35") "0.0₃0895"
Set-TextValue $ws.Range("E35") "  -6.17%  "

Set-TextValue $ws.Range("D36") "61.85"
Set-TextValue $ws.Range("E36") "  +4.56%  "

Set-TextValue $ws.Range("E37") "  +6.74%  "

Set-TextValue $ws.Range("D38") "40.65"
Set-TextValue $ws.Range("E38") "  -2.53%  "

Set-TextValue $ws.Range("E39") "  -0.47%  "

Set-TextValue $ws.Range("D40") "0.998"
Set-TextValue $ws.Range("E40") "  -0.06%  "

Set-TextValue $ws.Range("E41") "  +0.09%  "

Set-TextValue $ws.Range("E42") "  +2.27%  "

Set-TextValue $ws.Range("E43") "  +2.84%  "

Set-TextValue $ws.Range("E44") "  -2.30%  "

Set-TextValue $ws.Range("D45") "3.03"
Set-TextValue $ws.Range("E45") "  +2.58%  "

Set-TextValue $ws.Range("D46") "0.143"
Set-TextValue $ws.Range("E46") "  +0.20%  "

Set-TextValue $ws.Range("D47") "3.34"
Set-TextValue $ws.Range("E47") "  +5.01%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D48") "3.03"
Set-TextValue $ws.Range("E48") "  +6.29%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D49") "0.0₆0355"
Set-TextValue $ws.Range("E49") "  +9.11%  "

Set-TextValue $ws.Range("E50") "  -1.18%  "

Set-TextValue $ws.Range("E51") "  -3.25%  "
